# Auto-generated edit script: scheduled market-data refresh update
# Updates currentAveragePrice / LevePriceNQ / LevePriceHQ / LeveProfit* columns (H-N)
# across all 8 job sheets, per the upstream scraper run.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 251818.25
$ws.Range("I28").Value = 333757.84
$ws.Range("J28").Value = 5999.5
$ws.Range("K28").Value = 333757.84
$ws.Range("L28").Value = 5999.5
$ws.Range("M28").Value = -333272.84
$ws.Range("N28").Value = -6969.5
# Row 43
$ws.Range("H43").Value = 18458
$ws.Range("I43").Value = 16000
$ws.Range("J43").Value = 19277.334
$ws.Range("K43").Value = 16000
$ws.Range("L43").Value = 19277.334
$ws.Range("M43").Value = -15931
$ws.Range("N43").Value = -19415.334
# Row 132
$ws.Range("H132").Value = 2150.1667
$ws.Range("I132").Value = 1827.3243
$ws.Range("K132").Value = 5481.9729
$ws.Range("M132").Value = -2951.9729
# Row 137
$ws.Range("H137").Value = 37773.305
$ws.Range("I137").Value = 64872.69
$ws.Range("J137").Value = 2544.1
$ws.Range("K137").Value = 194618.07
$ws.Range("L137").Value = 7632.299999999999
$ws.Range("M137").Value = -192068.07
$ws.Range("N137").Value = -12732.3
# Row 138
$ws.Range("H138").Value = 4139.1
$ws.Range("I138").Value = 2696.35
$ws.Range("J138").Value = 4499.7876
$ws.Range("K138").Value = 8089.049999999999
$ws.Range("L138").Value = 13499.3628
$ws.Range("M138").Value = -2949.049999999999
$ws.Range("N138").Value = -23779.3628

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 18299690
$ws.Range("I32").Value = 18146012
$ws.Range("J32").Value = 21168334
$ws.Range("K32").Value = 18146012
$ws.Range("L32").Value = 21168334
$ws.Range("M32").Value = -18145725
$ws.Range("N32").Value = -21168908
# Row 45
$ws.Range("H45").Value = 4643.2915
$ws.Range("I45").Value = 4419.4736
$ws.Range("J45").Value = 5493.8
$ws.Range("K45").Value = 4419.4736
$ws.Range("L45").Value = 5493.8
$ws.Range("M45").Value = -4042.4736
$ws.Range("N45").Value = -6247.8
# Row 63
$ws.Range("H63").Value = 5233.6
$ws.Range("I63").Value = 2497.5
$ws.Range("K63").Value = 2497.5
$ws.Range("M63").Value = -1811.5
# Row 66
$ws.Range("H66").Value = 5233.6
$ws.Range("I66").Value = 2497.5
$ws.Range("K66").Value = 12487.5
$ws.Range("M66").Value = -9055.5
# Row 88
$ws.Range("H88").Value = 3519.3333
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3519.3333
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 3519.3333
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -4331.3333
# Row 91
$ws.Range("H91").Value = 3519.3333
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3519.3333
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 3519.3333
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -6327.3333
# Row 102
$ws.Range("H102").Value = 2516.25
$ws.Range("I102").Value = 1923.8
$ws.Range("J102").Value = 3503.6667
$ws.Range("K102").Value = 1923.8
$ws.Range("L102").Value = 3503.6667
$ws.Range("M102").Value = -301.8
$ws.Range("N102").Value = -6747.6667

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 111.166664
$ws.Range("I11").Value = 111.166664
$ws.Range("K11").Value = 111.166664
$ws.Range("M11").Value = 28.833336
# Row 130
$ws.Range("H130").Value = 74784.75
$ws.Range("J130").Value = 74784.75
$ws.Range("L130").Value = 74784.75
$ws.Range("N130").Value = -84824.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 8492.25
$ws.Range("I122").Value = 4989
$ws.Range("K122").Value = 14967
$ws.Range("M122").Value = -12517
# Row 127
$ws.Range("H127").Value = 119971.75
$ws.Range("J127").Value = 119971.75
$ws.Range("L127").Value = 119971.75
$ws.Range("N127").Value = -129891.75
# Row 132
$ws.Range("H132").Value = 4130.108
$ws.Range("I132").Value = 3870.4
$ws.Range("J132").Value = 5243.143
$ws.Range("K132").Value = 11611.2
$ws.Range("L132").Value = 15729.429
$ws.Range("M132").Value = -9081.200000000001
$ws.Range("N132").Value = -20789.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 51
$ws.Range("H51").Value = 1967.8334
$ws.Range("I51").Value = 935.6667
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 2807.0001
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -2347.0001
$ws.Range("N51").Value = -9920
# Row 68
$ws.Range("H68").Value = 2957.375
$ws.Range("J68").Value = 5333
$ws.Range("L68").Value = 15999
$ws.Range("N68").Value = -17621
# Row 71
$ws.Range("H71").Value = 2957.375
$ws.Range("J71").Value = 5333
$ws.Range("L71").Value = 47997
$ws.Range("N71").Value = -56109
# Row 107
$ws.Range("H107").Value = 1399.96
$ws.Range("I107").Value = 1106.3
$ws.Range("J107").Value = 1595.7333
$ws.Range("K107").Value = 3318.9
$ws.Range("L107").Value = 4787.199900000001
$ws.Range("M107").Value = -1398.9
$ws.Range("N107").Value = -8627.1999
# Row 139
$ws.Range("H139").Value = 2602.8262
$ws.Range("I139").Value = 1943.25
$ws.Range("K139").Value = 5829.75
$ws.Range("M139").Value = -689.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 36929.61
$ws.Range("I113").Value = 4513.6665
$ws.Range("K113").Value = 4513.6665
$ws.Range("M113").Value = -2343.6665
# Row 122
$ws.Range("H122").Value = 1899
$ws.Range("I122").Value = 1899
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5697
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3247
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 3315.1667
$ws.Range("I126").Value = 3231.889
$ws.Range("K126").Value = 9695.667000000001
$ws.Range("M126").Value = -7225.667000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1943.7727
$ws.Range("J22").Value = 2351
$ws.Range("L22").Value = 2351
$ws.Range("N22").Value = -2941
# Row 27
$ws.Range("H27").Value = 1943.7727
$ws.Range("J27").Value = 2351
$ws.Range("L27").Value = 2351
$ws.Range("N27").Value = -2565
# Row 55
$ws.Range("H55").Value = 1146.8148
$ws.Range("I55").Value = 1155.8948
$ws.Range("J55").Value = 1125.25
$ws.Range("K55").Value = 1155.8948
$ws.Range("L55").Value = 1125.25
$ws.Range("M55").Value = -982.8948
$ws.Range("N55").Value = -1471.25
# Row 94
$ws.Range("H94").Value = 99500
$ws.Range("J94").Value = 99500
$ws.Range("L94").Value = 99500
$ws.Range("N94").Value = -100852
# Row 122
$ws.Range("H122").Value = 19136.727
$ws.Range("I122").Value = 29000.6
$ws.Range("J122").Value = 10916.833
$ws.Range("K122").Value = 87001.79999999999
$ws.Range("L122").Value = 32750.499
$ws.Range("M122").Value = -84551.79999999999
$ws.Range("N122").Value = -37650.499

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 95360.63
$ws.Range("J81").Value = 5142.857
$ws.Range("L81").Value = 10285.714
$ws.Range("N81").Value = -12407.714
# Row 84
$ws.Range("H84").Value = 95360.63
$ws.Range("J84").Value = 5142.857
$ws.Range("L84").Value = 51428.57
$ws.Range("N84").Value = -62036.57
# Row 113
$ws.Range("H113").Value = 901.3333
$ws.Range("I113").Value = 901.3333
$ws.Range("K113").Value = 2703.9999
$ws.Range("M113").Value = -533.9998999999998
# Row 130
$ws.Range("H130").Value = 54142.668
$ws.Range("J130").Value = 54142.668
$ws.Range("L130").Value = 54142.668
$ws.Range("N130").Value = -64182.668
# Row 136
$ws.Range("H136").Value = 37416.066
$ws.Range("I136").Value = 3783.158
$ws.Range("K136").Value = 11349.474
$ws.Range("M136").Value = -8799.474
